$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column; existing A:D shift right to B:E
$ws.Columns("A:A").Insert()

# New narrow "Tab name" column
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Updated (rewritten) Cypher queries for the Cases tab / Stat query cells
# (now living in B2/C2 after the column insert above)
$casesQuery = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "ASIAN"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@

$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "ASIAN"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@

$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

# Row 2 grew taller to fit the longer wrapped query text
$ws.Rows("2:2").RowHeight = 174

# New, narrow best-fit width for the TabName column
$ws.Columns("A:A").ColumnWidth = 8

# Selection moves to the (now relocated) first big query cell
$ws.Range("B2").Select()
